# Applies updated cryptocurrency price/volume figures to the "cryptos" sheet.
# Source cells are plain text (inlineStr) cells, e.g. "330.95", "-0.06%".
# A leading apostrophe forces Excel/COM to keep the new value as literal text
# (matching the original cell type) instead of re-interpreting it as a number
# or percentage, which would change the stored value/format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.83"
$ws.Range("E2").Value = "'0.19%"
$ws.Range("D3").Value = "'41.45"
$ws.Range("E3").Value = "'0.35%"
$ws.Range("D4").Value = "'5.662"
$ws.Range("E4").Value = "'-1.59%"
$ws.Range("D5").Value = "'0.08344"
$ws.Range("E5").Value = "'2.97%"
$ws.Range("D6").Value = "'8.776"
$ws.Range("E6").Value = "'0.47%"
$ws.Range("D7").Value = "'2.003"
$ws.Range("E7").Value = "'-3.32%"
$ws.Range("D8").Value = "'4.503"
$ws.Range("E8").Value = "'-0.33%"
$ws.Range("D9").Value = "'2.943"
$ws.Range("E9").Value = "'-0.37%"
$ws.Range("D10").Value = "'0.9255"
$ws.Range("E10").Value = "'0.02%"
$ws.Range("D11").Value = "'0.1294"
$ws.Range("E11").Value = "'0.90%"
$ws.Range("D12").Value = "'0.1960"
$ws.Range("E12").Value = "'0.15%"
$ws.Range("D13").Value = "'0.09350"
$ws.Range("E13").Value = "'1.76%"
$ws.Range("D14").Value = "'0.03884"
$ws.Range("E14").Value = "'4.98%"
$ws.Range("D15").Value = "'0.1059"
$ws.Range("E15").Value = "'0.82%"
$ws.Range("E16").Value = "'-0.43%"
$ws.Range("D17").Value = "'0.006064"
$ws.Range("E17").Value = "'-1.85%"
$ws.Range("D18").Value = "'3.443"
$ws.Range("E18").Value = "'1.71%"
$ws.Range("E19").Value = "'0.23%"
$ws.Range("D20").Value = "'8.545"
$ws.Range("E20").Value = "'-3.12%"
$ws.Range("D21").Value = "'0.1354"
$ws.Range("E21").Value = "'-1.75%"
$ws.Range("E22").Value = "'-6.30%"
$ws.Range("D23").Value = "'0.04416"
$ws.Range("E23").Value = "'-0.42%"
$ws.Range("E24").Value = "'1.36%"
$ws.Range("D25").Value = "'0.004384"
$ws.Range("E25").Value = "'-1.86%"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'-3.27%"
$ws.Range("D39").Value = "'0.02802"
$ws.Range("E39").Value = "'0.38%"
$ws.Range("D40").Value = "'0.05527"
$ws.Range("E40").Value = "'0.11%"
$ws.Range("D41").Value = "'0.007805"
$ws.Range("E41").Value = "'1.88%"
$ws.Range("E42").Value = "'1.14%"
$ws.Range("D43").Value = "'0.009310"
$ws.Range("E43").Value = "'-5.67%"
$ws.Range("D44").Value = "'0.002071"
$ws.Range("E44").Value = "'-5.18%"
$ws.Range("D45").Value = "'0.01108"
$ws.Range("E45").Value = "'-6.92%"
$ws.Range("D46").Value = "'0.00006986"
$ws.Range("E46").Value = "'3.04%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.18%"
$ws.Range("E48").Value = "'9.19%"
$ws.Range("D49").Value = "'0.002278"
$ws.Range("E49").Value = "'0.10%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.16%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.18%"
